$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Clear A2 (was 1951) and B2 (was "Quách Bảo Hưng")
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()

# Update IdCardNo-like value in E2
$ws.Range("E2").Value = 46200608014

# Update AM2 and AN2
$ws.Range("AM2").Value = 1
$ws.Range("AN2").Value = 0

# Clear AQ2 (was "Thẻ BHYT hợp lệ")
$ws.Range("AQ2").ClearContents()

# Update sheet view selection to AN2 (matches the saved workbook state)
$ws.Activate()
$ws.Range("AN2").Select()
